$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rate matrix")

# ---------------------------------------------------------------------------
# 1. Structural change: insert a new compartment "asymptomatic" between
#    "exposed" and "symptomatic" -> new column before D, new row before 4.
# ---------------------------------------------------------------------------
$ws.Columns("D:D").Insert()
$ws.Rows("4:4").Insert()

# The insert shifted the old (unused, style-only) column J and L two columns
# to the right, leaving orphaned placeholder cells in K:M. Remove them - the
# new layout only spans A:J.
$ws.Columns("K:M").Delete()

# ---------------------------------------------------------------------------
# 2. Header row: label the new column D "asymptomatic"
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "asymptomatic"

# ---------------------------------------------------------------------------
# 3. Row 3 ("exposed") - rework formulas now that column D holds real data
# ---------------------------------------------------------------------------
$ws.Range("B3").ClearContents()
$ws.Range("B3").ClearFormats()

$ws.Range("C3").Formula = "=-SUM(D3:J3)"

$ws.Range("D3").Formula = "=1/8.5"
$ws.Range("D3").ClearFormats()

# E3 inherited the old D3 formula ("=1/11") when the column was inserted -
# that rate now belongs to the new D column, so E3 goes back to being blank.
$ws.Range("E3").ClearContents()

$ws.Range("J3").ClearFormats()
$ws.Range("J3").Formula = "=(0.8*(1/(68.5*52*7)))+(0.2*0.02*1/22)"
$ws.Range("J3").NumberFormat = "0.000000"

# ---------------------------------------------------------------------------
# 4. Row 4 (new "asymptomatic" compartment row)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "asymptomatic"
$ws.Range("B4").Value = 1.2
$ws.Range("B4").ClearFormats()

$ws.Range("C4").ClearContents()
$ws.Range("C4").ClearFormats()

$ws.Range("D4").Formula = "=-SUM(E4:J4)"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Formula = "=1/2.5"
$ws.Range("E4").ClearFormats()

$ws.Range("F4:I4").ClearContents()
$ws.Range("F4:I4").ClearFormats()

$ws.Range("J4").ClearFormats()
$ws.Range("J4").Formula = "=(0.8*(1/(68.5*52*7)))+(0.2*0.02*1/22)"
$ws.Range("J4").NumberFormat = "0.000000"

# Column D was inserted with the formatting of column C copied down the
# entire sheet (style applied to every row) - only rows 3 and 4 actually
# hold data now, so strip the inherited style from the remaining rows.
$ws.Range("D5:D9").ClearFormats()

# ---------------------------------------------------------------------------
# 5. Row 5 ("symptomatic", formerly row 4) - inline the literal rate instead
#    of referencing the old I8/J9 cell
# ---------------------------------------------------------------------------
$ws.Range("J5").Formula = "=(0.8*(1/(68.5*52*7)))+(0.2*0.02*1/22)"

# ---------------------------------------------------------------------------
# 6. Selection / view state to match the saved workbook
# ---------------------------------------------------------------------------
$ws.Range("H16").Select()

Write-Host "edit applied"
